$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Lingo")

# ---------------------------------------------------------------------------
# 1) Propagate formatting to the rows that are about to be created so the
#    new cells pick up the same visual style as their closest existing
#    "template" row (done first, before any values are overwritten).
# ---------------------------------------------------------------------------

# Rows 7-9 (Italy) reuse the exact formatting pattern that row 4 (the old
# "Germany" row) had (its D column uses the alternate percentage style).
$ws.Range("A4:K4").Copy()
$ws.Range("A7:K9").PasteSpecial(-4122)   # xlPasteFormats

# Rows 4-6 (Australia / extra Austria rows) reuse the plain formatting
# pattern used by row 1 (non-bold, "normal" D column style).
$ws.Range("A1:K1").Copy()
$ws.Range("A4:K6").PasteSpecial(-4122)   # xlPasteFormats

# Row 3 loses the bold country-name styling it previously had (that data
# belonged to Estonia, which is being replaced) so it matches the plain
# style used everywhere else in the refreshed table.
$ws.Range("B1:C1").Copy()
$ws.Range("B3:C3").PasteSpecial(-4122)   # xlPasteFormats

$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# 2) Overwrite the cell values for rows 2-9 with the refreshed Lingo table.
#    Row 1 (Albania / Manufacturer) is untouched.
# ---------------------------------------------------------------------------

$data = @(
    @(4,   "Albania",   "ALB", 0.15, 4280,  "Distributor",  2, "ALB2", 0.03004665490544474,  0.03939890539485347,   0.04903996252498364),
    @(6,   "Albania",   "ALB", 0.15, 4280,  "Principal",    3, "ALB3", 0.03177525212101811,  0.07772686245920965,   0.17564667785103982),
    @(27,  "Australia", "AUS", 0.3,  60050, "Principal",    3, "AUS3", 0.04840701411739386,  0.1278324518399418,    0.23334164580493164),
    @(28,  "Austria",   "AUT", 0.25, 47260, "Distributor",  2, "AUT2", 0.007677666919615235, 0.02981523327329755,   0.06295094632553316),
    @(29,  "Austria",   "AUT", 0.25, 47260, "Manufacturer", 1, "AUT1", 0.08171473421845657,  0.11700606768122054,   0.15229740114398455),
    @(193, "Italy",     "ITA", 0.279, 32830, "Distributor",  2, "ITA2", 0.00959587404677274, 0.02000235955350866,   0.04066699661131935),
    @(194, "Italy",     "ITA", 0.279, 32830, "Manufacturer", 1, "ITA1", 0.03972062576059085, 0.061999169361979495,  0.10827656744114039),
    @(195, "Italy",     "ITA", 0.279, 32830, "Principal",    3, "ITA3", 0.02191384689688711, 0.048206183032513075,  0.0997216190843369)
)

$r = 2
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]
    $ws.Cells.Item($r, 8).Value = $row[7]
    $ws.Cells.Item($r, 9).Value = $row[8]
    $ws.Cells.Item($r, 10).Value = $row[9]
    $ws.Cells.Item($r, 11).Value = $row[10]
    $r = $r + 1
}

# ---------------------------------------------------------------------------
# 3) Update the defined names that referenced the old Lingo!$*$1:$*$4 ranges
#    so they now span the grown table (rows 1-9).
# ---------------------------------------------------------------------------

$wb.Names.Item("ct").RefersTo  = "=Lingo!`$D`$1:`$D`$9"
$wb.Names.Item("gni").RefersTo = "=Lingo!`$E`$1:`$E`$9"
$wb.Names.Item("pli").RefersTo = "=Lingo!`$G`$1:`$G`$9"
$wb.Names.Item("pl").RefersTo  = "=Lingo!`$I`$1:`$I`$9"
$wb.Names.Item("pm").RefersTo  = "=Lingo!`$J`$1:`$J`$9"
$wb.Names.Item("pu").RefersTo  = "=Lingo!`$K`$1:`$K`$9"

# ---------------------------------------------------------------------------
# 4) Update sheet view / selection state to match the edited workbook:
#    - "Final (2)" scrolled down with a multi-row selection.
#    - "Lingo" becomes the active sheet, scrolled to/selecting C11.
# ---------------------------------------------------------------------------

$wsFinal = $wb.Worksheets.Item("Final (2)")
$wsFinal.Activate()
$wsFinal.Range("A194:XFD196").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 155

$ws.Activate()
$ws.Range("C11").Select() | Out-Null
